$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "NX-OS-21"
$ws.Range("B23").Value = "GigabitEthernet 0/0/20"
$ws.Range("C23").Value = 41
$ws.Range("D23").Value = "DOWN"

$ws.Range("A23:D23").Style = $ws.Range("A22:D22").Style

$ws.Range("G24").Select()
